$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# Insert a new keyword row "Seafo" / "Food & Drink" right before the existing
# "Sweet" row (worksheet row 43), shifting all subsequent rows down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(43).Insert()
$lo.Resize($ws.Range("A1:C90"))
$ws.Range("A43").Value2 = "Seafo"
$ws.Range("C43").Value2 = "Food & Drink"
$ws.Range("B43").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"

# ---------------------------------------------------------------------------
# Insert a new keyword row "Priso" / "Government" right after the existing
# "Polic" row (now worksheet row 50, after the previous insert), shifting all
# subsequent rows down by one more.
# ---------------------------------------------------------------------------
$ws.Rows.Item(51).Insert()
$lo.Resize($ws.Range("A1:C91"))
$ws.Range("A51").Value2 = "Priso"
$ws.Range("C51").Value2 = "Government"
$ws.Range("B51").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"

# ---------------------------------------------------------------------------
# Conditional formatting ranges are anchored to explicit row numbers and do
# not auto-shift when rows are inserted, so re-point each rule's range to the
# post-insert location (mirrors what Excel does automatically on a real
# row-insert).
# ---------------------------------------------------------------------------
$fcs = $ws.Cells.FormatConditions
$fcs.Item(3).ModifyAppliesToRange($ws.Range("B71:B91,B51:B69,B2:B49"))
$fcs.Item(4).ModifyAppliesToRange($ws.Range("A89:A91,A78:A86,A31:A35,A38:A49,A71:A76,A2:A18,A20:A28,A52:A69"))
$fcs.Item(5).ModifyAppliesToRange($ws.Range("B70"))
$fcs.Item(6).ModifyAppliesToRange($ws.Range("B50"))
$fcs.Item(7).ModifyAppliesToRange($ws.Range("A50:A51"))

# ---------------------------------------------------------------------------
# Restore the view state: scrolled down a bit with A44 selected.
# ---------------------------------------------------------------------------
$ws.Range("A44").Select()
